$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price (D) column to text format first so values like "27.711.93"
# and trailing-zero numbers such as "153.40" are preserved verbatim,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.711.93'
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').Value = '1.905.10'
$ws.Range('E3').Value = '  +0.53%  '

$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').Value = '312.50'
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.14%  '

$ws.Range('D7').Value = '0.5201'
$ws.Range('E7').Value = '  +7.39%  '

$ws.Range('D8').Value = '0.3778'
$ws.Range('E8').Value = '  -0.45%  '

$ws.Range('D9').Value = '0.07241'
$ws.Range('E9').Value = '  -1.30%  '

$ws.Range('D10').Value = '21.33'
$ws.Range('E10').Value = '  +3.88%  '

$ws.Range('D11').Value = '0.9005'
$ws.Range('E11').Value = '  -1.64%  '

$ws.Range('D12').Value = '0.07628'
$ws.Range('E12').Value = '  -0.78%  '

$ws.Range('D13').Value = '1.896.26'
$ws.Range('E13').Value = '  +0.38%  '

$ws.Range('D14').Value = '5.447'
$ws.Range('E14').Value = '  -0.39%  '

$ws.Range('D15').Value = '92.06'
$ws.Range('E15').Value = '  +1.25%  '

$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  -0.15%  '

$ws.Range('D17').Value = '0.000008697'
$ws.Range('E17').Value = '  -1.25%  '

$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').Value = '27.748.04'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('D20').Value = '14.48'
$ws.Range('E20').Value = '  +0.13%  '

$ws.Range('D21').Value = '5.143'
$ws.Range('E21').Value = '  +0.51%  '

$ws.Range('D22').Value = '2.132.29'
$ws.Range('E22').Value = '  +0.90%  '

$ws.Range('D23').Value = '10.83'
$ws.Range('E23').Value = '  +0.60%  '

$ws.Range('D24').Value = '6.582'
$ws.Range('E24').Value = '  -0.21%  '

$ws.Range('D25').Value = '153.40'
$ws.Range('E25').Value = '  -0.33%  '

$ws.Range('D26').Value = '1.883'
$ws.Range('E26').Value = '  -1.07%  '

$ws.Range('D27').Value = '18.31'
$ws.Range('E27').Value = '  -0.30%  '

$ws.Range('D28').Value = '2.163'
$ws.Range('E28').Value = '  +1.38%  '

$ws.Range('D29').Value = '114.47'
$ws.Range('E29').Value = '  -1.17%  '

$ws.Range('D30').Value = '4.852'
$ws.Range('E30').Value = '  -0.90%  '

$ws.Range('D31').Value = '0.08981'
$ws.Range('E31').Value = '  +0.74%  '

$ws.Range('D32').Value = '4.859'
$ws.Range('E32').Value = '  +4.58%  '

$ws.Range('D33').Value = '3.172'
$ws.Range('E33').Value = '  +0.68%  '

$ws.Range('D34').Value = '1.230'
$ws.Range('E34').Value = '  +0.22%  '

$ws.Range('D35').Value = '0.7705'
$ws.Range('E35').Value = '  +0.83%  '

$ws.Range('D36').Value = '2.640'
$ws.Range('E36').Value = '  +4.47%  '

$ws.Range('D37').Value = '0.02082'
$ws.Range('E37').Value = '  +1.91%  '

$ws.Range('D38').Value = '3.061'
$ws.Range('E38').Value = '  +2.39%  '

$ws.Range('D39').Value = '1.094'
$ws.Range('E39').Value = '  -0.06%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5514'
$ws.Range('E40').Value = '  +0.82%  '

$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.05281'
$ws.Range('E41').Value = '  +0.28%  '

$ws.Range('D42').Value = '6.663'
$ws.Range('E42').Value = '  -3.70%  '

$ws.Range('D43').Value = '114.56'
$ws.Range('E43').Value = '  +4.16%  '

$ws.Range('D44').Value = '8.518'
$ws.Range('E44').Value = '  +0.57%  '

$ws.Range('D45').Value = '0.1511'
$ws.Range('E45').Value = '  -0.41%  '

$ws.Range('D46').Value = '0.4801'
$ws.Range('E46').Value = '  +0.35%  '

$ws.Range('D47').Value = '10.42'
$ws.Range('E47').Value = '  -1.22%  '

$ws.Range('D48').Value = '0.9987'
$ws.Range('E48').Value = '  -0.14%  '

$ws.Range('D49').Value = '1.613'
$ws.Range('E49').Value = '  -1.31%  '

$ws.Range('D50').Value = '66.61'
$ws.Range('E50').Value = '  -1.06%  '

$ws.Range('D51').Value = '0.05987'
$ws.Range('E51').Value = '  -1.07%  '

# Restore the default "Normal" style on the Price column so no stray
# cell-level styling is left behind (matches original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"